# Update crypto price/volume figures per the latest scrape.
# Numeric-looking "Price" values are prefixed with a leading apostrophe so
# Excel stores them as plain text (matching the original inlineStr cells)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.835.01'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '1.637.27'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  -0.49%  '
$ws.Range('D5').Value = '''216.91'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('D9').Value = '''0.0623'
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('E10').Value = '  +3.49%  '
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '1.866.83'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '1.639.94'
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('D14').Value = '''4.12'
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('D17').Value = '26.834.17'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('D19').Value = '''219.03'
$ws.Range('E19').Value = '  +1.93%  '
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = '''6.75'
$ws.Range('E21').Value = '  +5.89%  '
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('D23').Value = '''2.43'
$ws.Range('E23').Value = '  +3.45%  '
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('D25').Value = '''147.06'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('E27').Value = '  +4.57%  '
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').Value = '''15.77'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('D30').Value = '''0.0503'
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('E32').Value = '  -1.39%  '
$ws.Range('E33').Value = '  +0.74%  '
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('D35').Value = '1.260.37'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('E37').Value = '  +2.15%  '
$ws.Range('D38').Value = '''0.533'
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').Value = '''0.807'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('E42').Value = '  +1.84%  '
$ws.Range('D43').Value = '1.782.48'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('D45').Value = '''61.62'
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('D46').Value = '''92.02'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('E48').Value = '  +14.05%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '''7.60'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('E51').Value = '  -0.29%  '
